# Updated cryptos list — refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active sheet, per the upstream data pull.
#
# A handful of the new Price strings (e.g. "610.00", "7.61", "0.0903") are
# syntactically valid numbers, and Excel's normal text->value coercion would
# silently convert a plain Range.Value assignment into a numeric cell. The
# source data models these Price cells as plain text (matching the sheet's
# existing inline-string cells), so for those specific values we prefix the
# literal with an apostrophe (Excel's "treat as text" quote-prefix) before
# assigning it, then reset the cell's Style back to "Normal" so the quote
# prefix doesn't leave a stray number-format/style behind on the cell.

$ws = $excel.ActiveWorkbook.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '67.844.28'; ForceText = $false },
    @{ Cell = "E2"; Value = '  +0.68%  '; ForceText = $false },
    @{ Cell = "D3"; Value = '3.518.22'; ForceText = $false },
    @{ Cell = "E3"; Value = '  +0.72%  '; ForceText = $false },
    @{ Cell = "E4"; Value = '  +0.07%  '; ForceText = $false },
    @{ Cell = "D5"; Value = '610.00'; ForceText = $true },
    @{ Cell = "E5"; Value = '  +1.34%  '; ForceText = $false },
    @{ Cell = "D6"; Value = '152.37'; ForceText = $true },
    @{ Cell = "E6"; Value = '  +2.37%  '; ForceText = $false },
    @{ Cell = "D7"; Value = '3.517.00'; ForceText = $false },
    @{ Cell = "E7"; Value = '  +0.69%  '; ForceText = $false },
    @{ Cell = "E8"; Value = '  +0.05%  '; ForceText = $false },
    @{ Cell = "E9"; Value = '  +1.72%  '; ForceText = $false },
    @{ Cell = "E10"; Value = '  +4.22%  '; ForceText = $false },
    @{ Cell = "D11"; Value = '7.61'; ForceText = $true },
    @{ Cell = "E11"; Value = '  +8.92%  '; ForceText = $false },
    @{ Cell = "E12"; Value = '  +3.00%  '; ForceText = $false },
    @{ Cell = "D13"; Value = '32.85'; ForceText = $true },
    @{ Cell = "E13"; Value = '  +4.61%  '; ForceText = $false },
    @{ Cell = "E14"; Value = '  +0.20%  '; ForceText = $false },
    @{ Cell = "D15"; Value = '4.116.89'; ForceText = $false },
    @{ Cell = "E15"; Value = '  +0.80%  '; ForceText = $false },
    @{ Cell = "D16"; Value = '3.524.78'; ForceText = $false },
    @{ Cell = "E16"; Value = '  +0.79%  '; ForceText = $false },
    @{ Cell = "D17"; Value = '67.957.32'; ForceText = $false },
    @{ Cell = "E17"; Value = '  +0.87%  '; ForceText = $false },
    @{ Cell = "D18"; Value = '0.116'; ForceText = $true },
    @{ Cell = "E18"; Value = '  -0.47%  '; ForceText = $false },
    @{ Cell = "D19"; Value = '6.60'; ForceText = $true },
    @{ Cell = "E19"; Value = '  +3.43%  '; ForceText = $false },
    @{ Cell = "D20"; Value = '15.63'; ForceText = $true },
    @{ Cell = "E20"; Value = '  +3.79%  '; ForceText = $false },
    @{ Cell = "D21"; Value = '9.81'; ForceText = $true },
    @{ Cell = "E21"; Value = '  +7.95%  '; ForceText = $false },
    @{ Cell = "D22"; Value = '450.47'; ForceText = $true },
    @{ Cell = "E22"; Value = '  +1.24%  '; ForceText = $false },
    @{ Cell = "D23"; Value = '0.635'; ForceText = $true },
    @{ Cell = "E23"; Value = '  +2.48%  '; ForceText = $false },
    @{ Cell = "D24"; Value = '78.34'; ForceText = $true },
    @{ Cell = "E24"; Value = '  +1.50%  '; ForceText = $false },
    @{ Cell = "E25"; Value = '  +2.42%  '; ForceText = $false },
    @{ Cell = "D26"; Value = '3.665.73'; ForceText = $false },
    @{ Cell = "E26"; Value = '  +0.90%  '; ForceText = $false },
    @{ Cell = "E27"; Value = '  -0.08%  '; ForceText = $false },
    @{ Cell = "D28"; Value = '9.00'; ForceText = $true },
    @{ Cell = "E28"; Value = '  +9.43%  '; ForceText = $false },
    @{ Cell = "D29"; Value = '10.14'; ForceText = $true },
    @{ Cell = "E29"; Value = '  +0.30%  '; ForceText = $false },
    @{ Cell = "E30"; Value = '  +9.59%  '; ForceText = $false },
    @{ Cell = "E31"; Value = '  +1.86%  '; ForceText = $false },
    @{ Cell = "D32"; Value = '0.170'; ForceText = $true },
    @{ Cell = "E32"; Value = '  +4.36%  '; ForceText = $false },
    @{ Cell = "E33"; Value = '  +0.02%  '; ForceText = $false },
    @{ Cell = "D34"; Value = '25.79'; ForceText = $true },
    @{ Cell = "E34"; Value = '  +0.90%  '; ForceText = $false },
    @{ Cell = "D35"; Value = '6.23'; ForceText = $true },
    @{ Cell = "E35"; Value = '  +3.05%  '; ForceText = $false },
    @{ Cell = "D36"; Value = '1.88'; ForceText = $true },
    @{ Cell = "E36"; Value = '  +2.97%  '; ForceText = $false },
    @{ Cell = "D37"; Value = '3.512.95'; ForceText = $false },
    @{ Cell = "E37"; Value = '  +1.08%  '; ForceText = $false },
    @{ Cell = "D38"; Value = '8.06'; ForceText = $true },
    @{ Cell = "E38"; Value = '  +0.28%  '; ForceText = $false },
    @{ Cell = "E39"; Value = '  +0.03%  '; ForceText = $false },
    @{ Cell = "D40"; Value = '2.33'; ForceText = $true },
    @{ Cell = "E40"; Value = '  +7.88%  '; ForceText = $false },
    @{ Cell = "E41"; Value = '  +0.01%  '; ForceText = $false },
    @{ Cell = "D42"; Value = '0.0903'; ForceText = $true },
    @{ Cell = "E42"; Value = '  +2.92%  '; ForceText = $false },
    @{ Cell = "D43"; Value = '173.43'; ForceText = $true },
    @{ Cell = "E43"; Value = '  -2.23%  '; ForceText = $false },
    @{ Cell = "D44"; Value = '5.55'; ForceText = $true },
    @{ Cell = "E44"; Value = '  +3.31%  '; ForceText = $false },
    @{ Cell = "D45"; Value = '30.65'; ForceText = $true },
    @{ Cell = "E45"; Value = '  +11.47%  '; ForceText = $false },
    @{ Cell = "D46"; Value = '0.883'; ForceText = $true },
    @{ Cell = "E46"; Value = '  +1.03%  '; ForceText = $false },
    @{ Cell = "D47"; Value = '47.43'; ForceText = $true },
    @{ Cell = "E47"; Value = '  +4.86%  '; ForceText = $false },
    @{ Cell = "D48"; Value = '1.31'; ForceText = $true },
    @{ Cell = "E48"; Value = '  +6.42%  '; ForceText = $false },
    @{ Cell = "D49"; Value = '2.55'; ForceText = $true },
    @{ Cell = "E49"; Value = '  +0.83%  '; ForceText = $false },
    @{ Cell = "D50"; Value = '7.69'; ForceText = $true },
    @{ Cell = "E50"; Value = '  +2.12%  '; ForceText = $false },
    @{ Cell = "D51"; Value = '0.256'; ForceText = $true },
    @{ Cell = "E51"; Value = '  +5.52%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
